$d = $word.ActiveDocument

# --- Edit the title / header block (paragraphs 1-4) ---
$p1 = $d.Paragraphs.Item(1)
$p1.Range.Text = 'Informe Técnico — Proyecto CodiGo'
$p1.Format.Alignment = 0  # wdAlignParagraphLeft -> drops the centered jc

$p2 = $d.Paragraphs.Item(2)
$p2.Range.Text = 'Autor: Juan Daniel Gómez Correa'

$p3 = $d.Paragraphs.Item(3)
$p3.Range.Text = 'Formación: Técnico en Programación de Software — SENA'

$p4 = $d.Paragraphs.Item(4)
$p4.Range.Text = 'Profesor: Juan Palo Jimenes Atehortúa'

# --- Drop everything from paragraph 5 onward (old body + images) ---
$delStart = $d.Paragraphs.Item(5).Range.Start
$delEnd = $d.Paragraphs.Item($d.Paragraphs.Count).Range.End
$d.Range($delStart, $delEnd).Delete()

# --- Rebuild the rest of the report from scratch ---
$d.Paragraphs.Last.Range.InsertParagraphAfter() | Out-Null
$cur = $d.Paragraphs.Last
$cur.Style = 'Normal'
$d.Paragraphs.Last.Range.InsertParagraphAfter() | Out-Null
$cur = $d.Paragraphs.Last
$cur.Style = 'Heading1'
$cur.Range.Text = '1. Estructura del Proyecto'
$d.Paragraphs.Last.Range.InsertParagraphAfter() | Out-Null
$cur = $d.Paragraphs.Last
$cur.Style = 'Normal'
$cur.Range.Text = '/index.html          → Página de inicio (login)'
$d.Paragraphs.Last.Range.InsertParagraphAfter() | Out-Null
$cur = $d.Paragraphs.Last
$cur.Style = 'Normal'
$cur.Range.Text = '/pages/register.html → Registro de usuario'
$d.Paragraphs.Last.Range.InsertParagraphAfter() | Out-Null
$cur = $d.Paragraphs.Last
$cur.Style = 'Normal'
$cur.Range.Text = '/pages/dashboard.html→ Panel principal (cursos)'
$d.Paragraphs.Last.Range.InsertParagraphAfter() | Out-Null
$cur = $d.Paragraphs.Last
$cur.Style = 'Normal'
$cur.Range.Text = '/pages/curso.html    → Vista de curso (teoría, consola y quiz)'
$d.Paragraphs.Last.Range.InsertParagraphAfter() | Out-Null
$cur = $d.Paragraphs.Last
$cur.Style = 'Normal'
$cur.Range.Text = '/assets/css/auth.css → Estilos para login/registro'
$d.Paragraphs.Last.Range.InsertParagraphAfter() | Out-Null
$cur = $d.Paragraphs.Last
$cur.Style = 'Normal'
$cur.Range.Text = '/assets/css/main.css → Estilos globales para dashboard y cursos'
$d.Paragraphs.Last.Range.InsertParagraphAfter() | Out-Null
$cur = $d.Paragraphs.Last
$cur.Style = 'Normal'
$cur.Range.Text = '/assets/js/login.js    → Lógica de login'
$d.Paragraphs.Last.Range.InsertParagraphAfter() | Out-Null
$cur = $d.Paragraphs.Last
$cur.Style = 'Normal'
$cur.Range.Text = '/assets/js/register.js → Lógica de registro'
$d.Paragraphs.Last.Range.InsertParagraphAfter() | Out-Null
$cur = $d.Paragraphs.Last
$cur.Style = 'Normal'
$cur.Range.Text = '/assets/js/dashboard.js→ Catálogo y progreso de cursos'
$d.Paragraphs.Last.Range.InsertParagraphAfter() | Out-Null
$cur = $d.Paragraphs.Last
$cur.Style = 'Normal'
$cur.Range.Text = '/assets/js/curso.js    → Teoría, consola y quiz del curso'
$d.Paragraphs.Last.Range.InsertParagraphAfter() | Out-Null
$cur = $d.Paragraphs.Last
$cur.Style = 'Normal'
$cur.Range.Text = 'Justificación del orden:'
$d.Paragraphs.Last.Range.InsertParagraphAfter() | Out-Null
$cur = $d.Paragraphs.Last
$cur.Style = 'Normal'
$cur.Range.Text = '- Separación de responsabilidades: cada vista tiene su propio HTML.'
$d.Paragraphs.Last.Range.InsertParagraphAfter() | Out-Null
$cur = $d.Paragraphs.Last
$cur.Style = 'Normal'
$cur.Range.Text = '- Carpetas `/assets/css` y `/assets/js` centralizan recursos para fácil mantenimiento.'
$d.Paragraphs.Last.Range.InsertParagraphAfter() | Out-Null
$cur = $d.Paragraphs.Last
$cur.Style = 'Normal'
$cur.Range.Text = '- CSS dividido en `auth.css` y `main.css`: uno exclusivo para autenticación y otro para interfaz principal.'
$d.Paragraphs.Last.Range.InsertParagraphAfter() | Out-Null
$cur = $d.Paragraphs.Last
$cur.Style = 'Normal'
$cur.Range.Text = '- JS modularizado: cada archivo gestiona solo la lógica de su pantalla.'
$d.Paragraphs.Last.Range.InsertParagraphAfter() | Out-Null
$cur = $d.Paragraphs.Last
$cur.Style = 'Heading1'
$cur.Range.Text = '2. Interfaces HTML'
$d.Paragraphs.Last.Range.InsertParagraphAfter() | Out-Null
$cur = $d.Paragraphs.Last
$cur.Style = 'Normal'
$cur.Range.Text = '2.1 index.html (Login)'
$d.Paragraphs.Last.Range.InsertParagraphAfter() | Out-Null
$cur = $d.Paragraphs.Last
$cur.Style = 'Normal'
$cur.Range.Text = '- Formulario con email y contraseña.'
$d.Paragraphs.Last.Range.InsertParagraphAfter() | Out-Null
$cur = $d.Paragraphs.Last
$cur.Style = 'Normal'
$cur.Range.Text = '- Validaciones visuales con mensajes de error.'
$d.Paragraphs.Last.Range.InsertParagraphAfter() | Out-Null
$cur = $d.Paragraphs.Last
$cur.Style = 'Normal'
$cur.Range.Text = '- Enlace a registro.'
$d.Paragraphs.Last.Range.InsertParagraphAfter() | Out-Null
$cur = $d.Paragraphs.Last
$cur.Style = 'Normal'
$cur.Range.Text = '2.2 register.html (Registro)'
$d.Paragraphs.Last.Range.InsertParagraphAfter() | Out-Null
$cur = $d.Paragraphs.Last
$cur.Style = 'Normal'
$cur.Range.Text = '- Campos: nombre, email, contraseña, confirmación.'
$d.Paragraphs.Last.Range.InsertParagraphAfter() | Out-Null
$cur = $d.Paragraphs.Last
$cur.Style = 'Normal'
$cur.Range.Text = '- Mensajes de error precisos.'
$d.Paragraphs.Last.Range.InsertParagraphAfter() | Out-Null
$cur = $d.Paragraphs.Last
$cur.Style = 'Normal'
$cur.Range.Text = '- Redirección al login tras registrarse.'
$d.Paragraphs.Last.Range.InsertParagraphAfter() | Out-Null
$cur = $d.Paragraphs.Last
$cur.Style = 'Normal'
$cur.Range.Text = '2.3 dashboard.html (Panel)'
$d.Paragraphs.Last.Range.InsertParagraphAfter() | Out-Null
$cur = $d.Paragraphs.Last
$cur.Style = 'Normal'
$cur.Range.Text = '- Encabezado fijo con botón de cerrar sesión.'
$d.Paragraphs.Last.Range.InsertParagraphAfter() | Out-Null
$cur = $d.Paragraphs.Last
$cur.Style = 'Normal'
$cur.Range.Text = '- Sección de cursos disponibles y sección de mis cursos.'
$d.Paragraphs.Last.Range.InsertParagraphAfter() | Out-Null
$cur = $d.Paragraphs.Last
$cur.Style = 'Normal'
$cur.Range.Text = '2.4 curso.html (Curso)'
$d.Paragraphs.Last.Range.InsertParagraphAfter() | Out-Null
$cur = $d.Paragraphs.Last
$cur.Style = 'Normal'
$cur.Range.Text = '- Encabezado con curso activo y nivel.'
$d.Paragraphs.Last.Range.InsertParagraphAfter() | Out-Null
$cur = $d.Paragraphs.Last
$cur.Style = 'Normal'
$cur.Range.Text = '- Bloques de Teoría, Consola virtual y Quiz.'
$d.Paragraphs.Last.Range.InsertParagraphAfter() | Out-Null
$cur = $d.Paragraphs.Last
$cur.Style = 'Heading1'
$cur.Range.Text = '3. Estilos (CSS)'
$d.Paragraphs.Last.Range.InsertParagraphAfter() | Out-Null
$cur = $d.Paragraphs.Last
$cur.Style = 'Normal'
$cur.Range.Text = '3.1 Paleta de colores — Teoría aplicada'
$d.Paragraphs.Last.Range.InsertParagraphAfter() | Out-Null
$cur = $d.Paragraphs.Last
$cur.Style = 'Normal'
$cur.Range.Text = '- Base oscura: concentración y reducción de fatiga visual.'
$d.Paragraphs.Last.Range.InsertParagraphAfter() | Out-Null
$cur = $d.Paragraphs.Last
$cur.Style = 'Normal'
$cur.Range.Text = '- Texto gris claro: accesible y de alto contraste.'
$d.Paragraphs.Last.Range.InsertParagraphAfter() | Out-Null
$cur = $d.Paragraphs.Last
$cur.Style = 'Normal'
$cur.Range.Text = '- Colores pastel vibrantes: verde menta (frescura y éxito), azul (confianza y tecnología), violeta (creatividad).'
$d.Paragraphs.Last.Range.InsertParagraphAfter() | Out-Null
$cur = $d.Paragraphs.Last
$cur.Style = 'Normal'
$cur.Range.Text = '- Rojo accesible para errores.'
$d.Paragraphs.Last.Range.InsertParagraphAfter() | Out-Null
$cur = $d.Paragraphs.Last
$cur.Style = 'Normal'
$cur.Range.Text = '- Uso de gradientes animados para dinamismo.'
$d.Paragraphs.Last.Range.InsertParagraphAfter() | Out-Null
$cur = $d.Paragraphs.Last
$cur.Style = 'Normal'
$cur.Range.Text = '3.2 Diseño visual'
$d.Paragraphs.Last.Range.InsertParagraphAfter() | Out-Null
$cur = $d.Paragraphs.Last
$cur.Style = 'Normal'
$cur.Range.Text = '- Botones primarios con gradientes animados.'
$d.Paragraphs.Last.Range.InsertParagraphAfter() | Out-Null
$cur = $d.Paragraphs.Last
$cur.Style = 'Normal'
$cur.Range.Text = '- Sombras suaves y bordes redondeados → estilo moderno amigable.'
$d.Paragraphs.Last.Range.InsertParagraphAfter() | Out-Null
$cur = $d.Paragraphs.Last
$cur.Style = 'Normal'
$cur.Range.Text = '- Responsive con media queries.'
$d.Paragraphs.Last.Range.InsertParagraphAfter() | Out-Null
$cur = $d.Paragraphs.Last
$cur.Style = 'Normal'
$cur.Range.Text = '- Consola virtual con estética de terminal real.'
$d.Paragraphs.Last.Range.InsertParagraphAfter() | Out-Null
$cur = $d.Paragraphs.Last
$cur.Style = 'Heading1'
$cur.Range.Text = '4. Lógica (JavaScript)'
$d.Paragraphs.Last.Range.InsertParagraphAfter() | Out-Null
$cur = $d.Paragraphs.Last
$cur.Style = 'Normal'
$cur.Range.Text = '4.1 login.js'
$d.Paragraphs.Last.Range.InsertParagraphAfter() | Out-Null
$cur = $d.Paragraphs.Last
$cur.Style = 'Normal'
$cur.Range.Text = '- Captura formulario y valida campos vacíos.'
$d.Paragraphs.Last.Range.InsertParagraphAfter() | Out-Null
$cur = $d.Paragraphs.Last
$cur.Style = 'Normal'
$cur.Range.Text = '- Comprueba credenciales contra localStorage.'
$d.Paragraphs.Last.Range.InsertParagraphAfter() | Out-Null
$cur = $d.Paragraphs.Last
$cur.Style = 'Normal'
$cur.Range.Text = '- Redirige al dashboard si es correcto.'
$d.Paragraphs.Last.Range.InsertParagraphAfter() | Out-Null
$cur = $d.Paragraphs.Last
$cur.Style = 'Normal'
$cur.Range.Text = '4.2 register.js'
$d.Paragraphs.Last.Range.InsertParagraphAfter() | Out-Null
$cur = $d.Paragraphs.Last
$cur.Style = 'Normal'
$cur.Range.Text = '- Valida nombre, email, contraseña fuerte y confirmación.'
$d.Paragraphs.Last.Range.InsertParagraphAfter() | Out-Null
$cur = $d.Paragraphs.Last
$cur.Style = 'Normal'
$cur.Range.Text = '- Evita correos duplicados.'
$d.Paragraphs.Last.Range.InsertParagraphAfter() | Out-Null
$cur = $d.Paragraphs.Last
$cur.Style = 'Normal'
$cur.Range.Text = '- Guarda en localStorage y redirige al login.'
$d.Paragraphs.Last.Range.InsertParagraphAfter() | Out-Null
$cur = $d.Paragraphs.Last
$cur.Style = 'Normal'
$cur.Range.Text = '4.3 dashboard.js'
$d.Paragraphs.Last.Range.InsertParagraphAfter() | Out-Null
$cur = $d.Paragraphs.Last
$cur.Style = 'Normal'
$cur.Range.Text = '- Maneja catálogo de cursos, disponibles y mis cursos.'
$d.Paragraphs.Last.Range.InsertParagraphAfter() | Out-Null
$cur = $d.Paragraphs.Last
$cur.Style = 'Normal'
$cur.Range.Text = '- Muestra progreso, XP y nivel.'
$d.Paragraphs.Last.Range.InsertParagraphAfter() | Out-Null
$cur = $d.Paragraphs.Last
$cur.Style = 'Normal'
$cur.Range.Text = '- Botones dinámicos: Inscribirme / Abrir / Cancelar.'
$d.Paragraphs.Last.Range.InsertParagraphAfter() | Out-Null
$cur = $d.Paragraphs.Last
$cur.Style = 'Normal'
$cur.Range.Text = '4.4 curso.js'
$d.Paragraphs.Last.Range.InsertParagraphAfter() | Out-Null
$cur = $d.Paragraphs.Last
$cur.Style = 'Normal'
$cur.Range.Text = '- Carga curso activo y bloquea contenido de otros.'
$d.Paragraphs.Last.Range.InsertParagraphAfter() | Out-Null
$cur = $d.Paragraphs.Last
$cur.Style = 'Normal'
$cur.Range.Text = '- Editor de código persistente con localStorage.'
$d.Paragraphs.Last.Range.InsertParagraphAfter() | Out-Null
$cur = $d.Paragraphs.Last
$cur.Style = 'Normal'
$cur.Range.Text = '- Consola aislada en iframe.'
$d.Paragraphs.Last.Range.InsertParagraphAfter() | Out-Null
$cur = $d.Paragraphs.Last
$cur.Style = 'Normal'
$cur.Range.Text = '- Quiz que otorga XP según respuestas correctas.'
$d.Paragraphs.Last.Range.InsertParagraphAfter() | Out-Null
$cur = $d.Paragraphs.Last
$cur.Style = 'Normal'
$cur.Range.Text = '- Sistema de niveles: Lv1 (0-99 XP), Lv2 (100-199 XP), Lv3 (200-240 XP).'
$d.Paragraphs.Last.Range.InsertParagraphAfter() | Out-Null
$cur = $d.Paragraphs.Last
$cur.Style = 'Heading1'
$cur.Range.Text = '5. Buenas Prácticas Aplicadas'
$d.Paragraphs.Last.Range.InsertParagraphAfter() | Out-Null
$cur = $d.Paragraphs.Last
$cur.Style = 'Normal'
$cur.Range.Text = '- Separación de capas (HTML, CSS, JS).'
$d.Paragraphs.Last.Range.InsertParagraphAfter() | Out-Null
$cur = $d.Paragraphs.Last
$cur.Style = 'Normal'
$cur.Range.Text = '- Accesibilidad: contrastes, feedback visual, placeholders claros.'
$d.Paragraphs.Last.Range.InsertParagraphAfter() | Out-Null
$cur = $d.Paragraphs.Last
$cur.Style = 'Normal'
$cur.Range.Text = '- Persistencia en navegador con localStorage.'
$d.Paragraphs.Last.Range.InsertParagraphAfter() | Out-Null
$cur = $d.Paragraphs.Last
$cur.Style = 'Normal'
$cur.Range.Text = '- Diseño responsivo y moderno.'
$d.Paragraphs.Last.Range.InsertParagraphAfter() | Out-Null
$cur = $d.Paragraphs.Last
$cur.Style = 'Normal'
$cur.Range.Text = '- Código modular (responsabilidad única por archivo).'
$d.Paragraphs.Last.Range.InsertParagraphAfter() | Out-Null
$cur = $d.Paragraphs.Last
$cur.Style = 'Normal'
$cur.Range.Text = '- Validaciones estrictas para formularios.'
$d.Paragraphs.Last.Range.InsertParagraphAfter() | Out-Null
$cur = $d.Paragraphs.Last
$cur.Style = 'Normal'
$cur.Range.Text = '- Experiencia de usuario optimizada con feedback inmediato, progreso y gamificación.'
$d.Paragraphs.Last.Range.InsertParagraphAfter() | Out-Null
$cur = $d.Paragraphs.Last
$cur.Style = 'Heading1'
$cur.Range.Text = '6. Conclusiones'
$d.Paragraphs.Last.Range.InsertParagraphAfter() | Out-Null
$cur = $d.Paragraphs.Last
$cur.Style = 'Normal'
$cur.Range.Text = 'El proyecto demuestra:'
$d.Paragraphs.Last.Range.InsertParagraphAfter() | Out-Null
$cur = $d.Paragraphs.Last
$cur.Style = 'Normal'
$cur.Range.Text = '- Dominio de HTML semántico, CSS moderno (flex, grid, variables, animaciones) y JavaScript estructurado.'
$d.Paragraphs.Last.Range.InsertParagraphAfter() | Out-Null
$cur = $d.Paragraphs.Last
$cur.Style = 'Normal'
$cur.Range.Text = '- Aplicación de principios de usabilidad, accesibilidad y diseño visual.'
$d.Paragraphs.Last.Range.InsertParagraphAfter() | Out-Null
$cur = $d.Paragraphs.Last
$cur.Style = 'Normal'
$cur.Range.Text = '- Uso de persistencia en cliente para un sistema educativo básico.'
$d.Paragraphs.Last.Range.InsertParagraphAfter() | Out-Null
$cur = $d.Paragraphs.Last
$cur.Style = 'Normal'
$cur.Range.Text = '- Implementación de gamificación con XP, niveles y progreso.'
$d.Paragraphs.Last.Range.InsertParagraphAfter() | Out-Null
$cur = $d.Paragraphs.Last
$cur.Style = 'Normal'
$cur.Range.Text = 'Conclusión: CodiGo es una plataforma educativa web con autenticación, gestión de usuarios, panel de cursos y ejecución práctica de código, todo en una interfaz clara, atractiva y moderna.'

Write-Output ("Final paragraph count: " + $d.Paragraphs.Count)
